# Apply the "change name exposures categories + correct error in naming of
# canton Geneva" edit to the annual_deaths workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the exposure categories in column E:
#    "O" -> "Over 75"  (rows 2-28, the over-75 block)
#    "U" -> "Under 75" (rows 29-55, the under-75 block)
$ws.Range("E2:E28").Value = "Over 75"
$ws.Range("E29:E55").Value = "Under 75"

# 2) Correct the canton name "Geneva" -> "Genève" (rows 10 and 37).
$ws.Range("A10").Value = "Genève"
$ws.Range("A37").Value = "Genève"

# 3) Update the view/selection state to match the saved workbook: no frozen
#    top-left scroll position, active cell A10 selected.
$ws.Range("A10").Select()
